# Generate Report for Handback
# Renames the two handback entries:
#   80eb6098-346a-450e-af2f-c3a99c7b83f9  ->  a9d38774-1895-4a92-b9d2-ea54ed7a5a70
#   ce68870e-f5bc-44bb-9eea-10b9c930e382  ->  ffff0067dcba-e716-43b1-84d2-08e377e05498
# and refreshes the generated xliff filenames / timestamps that go with the
# new handback run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- Overview
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.md"
$overview.Range("B2").Value = "e2e\a9d38774-1895-4a92-b9d2-ea54ed7a5a70.md"
$overview.Range("G2").Value = "2016-08-26 01:02:17"

$overview.Range("A3").Value = "ffff0067dcba-e716-43b1-84d2-08e377e05498.md"
$overview.Range("B3").Value = "e2e\ffff0067dcba-e716-43b1-84d2-08e377e05498.md"
$overview.Range("G3").Value = "2016-08-26 01:02:17"

# ------------------------------------------------------------------ zh-cn
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.md"
$zhcn.Range("G2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-26 01:02:12"
$zhcn.Range("I2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.md"
$zhcn.Range("J2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 01:02:28"

$zhcn.Range("A3").Value = "ffff0067dcba-e716-43b1-84d2-08e377e05498.md"
$zhcn.Range("G3").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-26 01:02:12"
$zhcn.Range("I3").Value = "ffff0067dcba-e716-43b1-84d2-08e377e05498.md"
$zhcn.Range("J3").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 01:02:28"

# ------------------------------------------------------------------ de-de
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.md"
$dede.Range("G2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.de-de.xlf"
$dede.Range("H2").Value = "2016-08-26 01:02:17"
$dede.Range("I2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.md"
$dede.Range("J2").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 01:02:34"

$dede.Range("A3").Value = "ffff0067dcba-e716-43b1-84d2-08e377e05498.md"
$dede.Range("G3").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.de-de.xlf"
$dede.Range("H3").Value = "2016-08-26 01:02:17"
$dede.Range("I3").Value = "ffff0067dcba-e716-43b1-84d2-08e377e05498.md"
$dede.Range("J3").Value = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70.11d518846421860cb7a79ac95b58e1c940a6b15d.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 01:02:34"
